$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 725; this shifts the previous rows 725-766 down to 729-770,
# keeping all of their data (including any per-row formatting) intact.
$ws.Rows("725:728").Insert()

# Populate the 4 newly-inserted rows (725-728) with new weekly price data for Pera
# varieties sold at "Vega Modelo de Temuco" on 2022-07-11 (date serial 44753).

# Row 725: Forelle / Primera
$ws.Cells.Item(725,1).Value2  = 10
$ws.Cells.Item(725,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(725,3).Value2  = "La Araucanía"
$ws.Cells.Item(725,4).Value2  = 44753
$ws.Cells.Item(725,5).Value2  = 9
$ws.Cells.Item(725,6).Value2  = "Fruta"
$ws.Cells.Item(725,7).Value2  = 100104
$ws.Cells.Item(725,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(725,9).Value2  = 100104005
$ws.Cells.Item(725,10).Value2 = "Pera"
$ws.Cells.Item(725,11).Value2 = "Forelle"
$ws.Cells.Item(725,12).Value2 = "Primera"
$ws.Cells.Item(725,13).Value2 = 110
$ws.Cells.Item(725,14).Value2 = 12000
$ws.Cells.Item(725,15).Value2 = 12000
$ws.Cells.Item(725,16).Value2 = 12000
$ws.Cells.Item(725,17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(725,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(725,19).Value2 = 667
$ws.Cells.Item(725,20).Value2 = 18

# Row 726: Packham's Triumph / Especial
$ws.Cells.Item(726,1).Value2  = 10
$ws.Cells.Item(726,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(726,3).Value2  = "La Araucanía"
$ws.Cells.Item(726,4).Value2  = 44753
$ws.Cells.Item(726,5).Value2  = 9
$ws.Cells.Item(726,6).Value2  = "Fruta"
$ws.Cells.Item(726,7).Value2  = 100104
$ws.Cells.Item(726,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(726,9).Value2  = 100104005
$ws.Cells.Item(726,10).Value2 = "Pera"
$ws.Cells.Item(726,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(726,12).Value2 = "Especial"
$ws.Cells.Item(726,13).Value2 = 185
$ws.Cells.Item(726,14).Value2 = 12000
$ws.Cells.Item(726,15).Value2 = 12000
$ws.Cells.Item(726,16).Value2 = 12000
$ws.Cells.Item(726,17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(726,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(726,19).Value2 = 667
$ws.Cells.Item(726,20).Value2 = 18

# Row 727: Packham's Triumph / Primera
$ws.Cells.Item(727,1).Value2  = 10
$ws.Cells.Item(727,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(727,3).Value2  = "La Araucanía"
$ws.Cells.Item(727,4).Value2  = 44753
$ws.Cells.Item(727,5).Value2  = 9
$ws.Cells.Item(727,6).Value2  = "Fruta"
$ws.Cells.Item(727,7).Value2  = 100104
$ws.Cells.Item(727,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(727,9).Value2  = 100104005
$ws.Cells.Item(727,10).Value2 = "Pera"
$ws.Cells.Item(727,11).Value2 = "Packham's Triumph"
$ws.Cells.Item(727,12).Value2 = "Primera"
$ws.Cells.Item(727,13).Value2 = 320
$ws.Cells.Item(727,14).Value2 = 10000
$ws.Cells.Item(727,15).Value2 = 12000
$ws.Cells.Item(727,16).Value2 = 10969
$ws.Cells.Item(727,17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(727,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(727,19).Value2 = 609
$ws.Cells.Item(727,20).Value2 = 18

# Row 728: Winter Nelis / Primera
$ws.Cells.Item(728,1).Value2  = 10
$ws.Cells.Item(728,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(728,3).Value2  = "La Araucanía"
$ws.Cells.Item(728,4).Value2  = 44753
$ws.Cells.Item(728,5).Value2  = 9
$ws.Cells.Item(728,6).Value2  = "Fruta"
$ws.Cells.Item(728,7).Value2  = 100104
$ws.Cells.Item(728,8).Value2  = "Frutos de pepita"
$ws.Cells.Item(728,9).Value2  = 100104005
$ws.Cells.Item(728,10).Value2 = "Pera"
$ws.Cells.Item(728,11).Value2 = "Winter Nelis"
$ws.Cells.Item(728,12).Value2 = "Primera"
$ws.Cells.Item(728,13).Value2 = 155
$ws.Cells.Item(728,14).Value2 = 12000
$ws.Cells.Item(728,15).Value2 = 12000
$ws.Cells.Item(728,16).Value2 = 12000
$ws.Cells.Item(728,17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(728,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(728,19).Value2 = 667
$ws.Cells.Item(728,20).Value2 = 18
